$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.232.79'
$ws.Range('D3').Value = '1.829.90'
$ws.Range('E3').Value = '  -0.74%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.09'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  -0.99%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6099'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -3.24%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$__style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07090'
$ws.Range('D8').Style = $__style
$ws.Range('E8').Value = '  -4.97%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2831'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  -2.67%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.99'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  -3.84%  '
$__style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07641'
$ws.Range('D11').Style = $__style
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').Value = '1.833.17'
$ws.Range('E12').Value = '  -0.62%  '
$__style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.822'
$ws.Range('D13').Style = $__style
$ws.Range('E13').Value = '  -3.37%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6380'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  -6.07%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009956'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').Value = '2.075.02'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('E17').Value = '  -2.92%  '
$__style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.990'
$ws.Range('D18').Style = $__style
$ws.Range('E18').Value = '  -4.66%  '
$ws.Range('D19').Value = '29.206.05'
$ws.Range('E19').Value = '  -0.71%  '
$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '230.67'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('E21').Value = '  -4.16%  '
$ws.Range('E22').Value = '  +0.11%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.037'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  -5.11%  '
$ws.Range('E24').Value = '  +0.09%  '
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.54'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  -1.82%  '
$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.085'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  -4.78%  '
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1295'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  -4.29%  '
$ws.Range('E28').Value = '  -3.83%  '
$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06738'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  +3.03%  '
$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.468'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  +1.82%  '
$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.460'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  -1.85%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.846'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  -5.32%  '
$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.821'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  -6.20%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.131'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -0.83%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.731'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  -5.97%  '
$__style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6582'
$ws.Range('D36').Style = $__style
$ws.Range('E36').Value = '  -5.67%  '
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.551'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').Value = '1.238.41'
$ws.Range('E38').Value = '  -0.82%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.755'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('E40').Value = '  -4.71%  '
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.611'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  -2.51%  '
$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9303'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  -0.16%  '
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('D43').Style = $__style
$ws.Range('D44').Value = '1.984.83'
$ws.Range('E44').Value = '  -0.76%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.57'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  -0.24%  '
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.74'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  -2.82%  '
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('E48').Value = '  -4.97%  '
$__style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.554'
$ws.Range('D49').Style = $__style
$ws.Range('E49').Value = '  -4.68%  '
$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1088'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  -5.07%  '
$__style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05574'
$ws.Range('D51').Style = $__style
$ws.Range('E51').Value = '  -1.69%  '
